$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new program rows (TU Chemnitz Advanced Manufacturing, TUM Computational Mechanics)
$ws.Range("A7").Value = "TU_CHEMNITZ_ADVANCED_MANUFACTURING"
$ws.Range("B7").Value = "Yes"

$ws.Range("A8").Value = "TUM_COMPUTATIONAL_MECHANICS"
$ws.Range("B8").Value = "Yes"

# Extend the data validation list range to cover the new rows
$ws.Range("B1:B8").Validation.Delete()
$ws.Range("B1:B8").Validation.Add(3, 1, 1, '"Yes,No"')

# Move selection to reflect the final saved state
$ws.Range("C6").Select()
